$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 13
$ws.Range("C13").Value = "[]"
$ws.Range("E13").Value = "[]"

# Row 16
$ws.Range("C16").Value = "[336]"
$ws.Range("E16").Value = "[336]"

# Row 18
$ws.Range("C18").Value = "[]"
$ws.Range("D18").Value = "[]"
$ws.Range("E18").Value = "[]"
$ws.Range("F18").Value = "[90]"

# Row 19
$ws.Range("C19").Value = "[585]"

# Row 20
$ws.Range("C20").Value = "[]"
$ws.Range("D20").Value = "[]"
$ws.Range("E20").Value = "[]"
$ws.Range("F20").Value = "[700]"

# Row 21
$ws.Range("C21").Value = "[]"
$ws.Range("D21").Value = "[]"
$ws.Range("E21").Value = "[]"
$ws.Range("F21").Value = "[260]"

# Row 22
$ws.Range("C22").Value = "[]"
$ws.Range("D22").Value = "[]"
$ws.Range("E22").Value = "[]"
$ws.Range("F22").Value = "[170]"

# Row 24
$ws.Range("C24").Value = "[]"
$ws.Range("D24").Value = "[]"
$ws.Range("E24").Value = "[]"
$ws.Range("F24").Value = "[340]"

# Row 26
$ws.Range("C26").Value = "[483]"

# Row 27
$ws.Range("C27").Value = "[375]"

# Row 28
$ws.Range("C28").Value = "[291]"

# Row 30
$ws.Range("C30").Value = "[150]"

# Row 31
$ws.Range("C31").Value = "[375]"
$ws.Range("D31").Value = "[300]"
$ws.Range("F31").Value = "[]"

# Row 33
$ws.Range("C33").Value = "[]"
$ws.Range("D33").Value = "[]"
$ws.Range("E33").Value = "[]"
$ws.Range("F33").Value = "[600]"

# Row 34
$ws.Range("C34").Value = "[468]"
$ws.Range("E34").Value = "[468]"

# Row 36
$ws.Range("C36").Value = "[213]"

# Row 37
$ws.Range("C37").Value = "[]"
$ws.Range("E37").Value = "[]"

# Row 38
$ws.Range("C38").Value = "[]"
$ws.Range("E38").Value = "[]"

# Row 39
$ws.Range("C39").Value = "[]"
$ws.Range("E39").Value = "[]"
